$d = $word.ActiveDocument

# Replace the text of the run that currently holds $oldText (inside paragraph
# $paraIndex) with $newText, while keeping any other runs in that paragraph
# (e.g. a leading empty <w:r/>) untouched and in their original position.
# $rPrXml is the literal run-properties markup (e.g. '<w:rPr><w:b/></w:rPr>')
# to re-apply to the new run, or an empty string if the run has no formatting.
function Set-RunText($paraIndex, $oldText, $newText, $rPrXml) {
    $full = $d.Paragraphs($paraIndex).Range
    $oldLen = $oldText.Length
    $textRunRange = $d.Range($full.Start, $full.Start + $oldLen)

    $escaped = $newText -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $runXml = '<w:r>' + $rPrXml + '<w:t>' + $escaped + '</w:t></w:r>'
    $xmlFrag = '<?xml version="1.0"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $runXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $textRunRange.InsertXML($xmlFrag)
}

# 1. Title heading (Heading1, paragraph 1)
Set-RunText 1 "Play Lumber Jack for Free - Exciting Forest-Themed Slot Game" "Play Lumber Jack Slot Free - Review & Game Overview" ''

# 2. "What we like" bullet list items
Set-RunText 39 "High-quality graphics and polished design" "High-quality design" ''
Set-RunText 40 "10 adjustable paylines with bets ranging from €0.01 to €50" "Pleasant soundtrack" ''
Set-RunText 41 "RTP of 96.07% and medium volatility" "Adjustable paylines and bets" ''
Set-RunText 42 "Thrilling free spins and gamble feature" "Bonus game with free spins" ''

# 3. "What we don't like" bullet list items
Set-RunText 44 "Limited range of multipliers during base spins" "Medium volatility" ''
Set-RunText 45 "The game may not appeal to casual players due to its medium volatility" "Limited multiplier symbols during base spins" ''

# 4. Bold title repeated near the end (paragraph 46)
Set-RunText 46 "Play Lumber Jack for Free - Exciting Forest-Themed Slot Game" "Play Lumber Jack Slot Free - Review & Game Overview" '<w:rPr><w:b/></w:rPr>'

# 5. Meta description italic paragraph (paragraph 47)
Set-RunText 47 "Read our review of Lumber Jack, a forest-themed slot game with high-quality design. Play Lumber Jack for free and enjoy free spins and a gamble feature." "Read our review of Lumber Jack slot game and play it for free. Discover its features and bonus game with free spins." '<w:rPr><w:i/></w:rPr>'
